$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 74
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(74, 3).Value = 'Bíobío'
$ws.Cells.Item(74, 4).Value = 44644
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112024
$ws.Cells.Item(74, 7).Value = 'Choclo'
$ws.Cells.Item(74, 8).Value = 'Choclero'
$ws.Cells.Item(74, 9).Value = 'Primera'
$ws.Cells.Item(74, 10).Value = 2500
$ws.Cells.Item(74, 11).Value = 200
$ws.Cells.Item(74, 12).Value = 250
$ws.Cells.Item(74, 13).Value = 230
$ws.Cells.Item(74, 14).Value = '$/unidad'
$ws.Cells.Item(74, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(74, 16).Value = 230
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = 'Hortaliza'

# Row 75
$ws.Cells.Item(75, 1).Value = 11
$ws.Cells.Item(75, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(75, 3).Value = 'Bíobío'
$ws.Cells.Item(75, 4).Value = 44355
$ws.Cells.Item(75, 5).Value = 8
$ws.Cells.Item(75, 6).Value = 100112024
$ws.Cells.Item(75, 7).Value = 'Choclo'
$ws.Cells.Item(75, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(75, 9).Value = 'Primera'
$ws.Cells.Item(75, 10).Value = 50
$ws.Cells.Item(75, 11).Value = 11000
$ws.Cells.Item(75, 12).Value = 12000
$ws.Cells.Item(75, 13).Value = 11400
$ws.Cells.Item(75, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(75, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(75, 16).Value = 190
$ws.Cells.Item(75, 17).Value = 60
$ws.Cells.Item(75, 18).Value = 'Hortaliza'

# Row 76
$ws.Cells.Item(76, 1).Value = 11
$ws.Cells.Item(76, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(76, 3).Value = 'Bíobío'
$ws.Cells.Item(76, 4).Value = 44483
$ws.Cells.Item(76, 5).Value = 8
$ws.Cells.Item(76, 6).Value = 100112024
$ws.Cells.Item(76, 7).Value = 'Choclo'
$ws.Cells.Item(76, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(76, 9).Value = 'Primera'
$ws.Cells.Item(76, 10).Value = 450
$ws.Cells.Item(76, 11).Value = 25000
$ws.Cells.Item(76, 12).Value = 26000
$ws.Cells.Item(76, 13).Value = 25556
$ws.Cells.Item(76, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(76, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(76, 16).Value = 365
$ws.Cells.Item(76, 17).Value = 70
$ws.Cells.Item(76, 18).Value = 'Hortaliza'

# Row 77
$ws.Cells.Item(77, 1).Value = 11
$ws.Cells.Item(77, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(77, 3).Value = 'Bíobío'
$ws.Cells.Item(77, 4).Value = 44617
$ws.Cells.Item(77, 5).Value = 8
$ws.Cells.Item(77, 6).Value = 100112024
$ws.Cells.Item(77, 7).Value = 'Choclo'
$ws.Cells.Item(77, 8).Value = 'Choclero'
$ws.Cells.Item(77, 9).Value = 'Primera'
$ws.Cells.Item(77, 10).Value = 20000
$ws.Cells.Item(77, 11).Value = 230
$ws.Cells.Item(77, 12).Value = 250
$ws.Cells.Item(77, 13).Value = 240
$ws.Cells.Item(77, 14).Value = '$/unidad'
$ws.Cells.Item(77, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(77, 16).Value = 240
$ws.Cells.Item(77, 17).Value = 1
$ws.Cells.Item(77, 18).Value = 'Hortaliza'

# Row 78
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(78, 3).Value = 'Bíobío'
$ws.Cells.Item(78, 4).Value = 44617
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = 100112024
$ws.Cells.Item(78, 7).Value = 'Choclo'
$ws.Cells.Item(78, 8).Value = 'Choclero'
$ws.Cells.Item(78, 9).Value = 'Segunda'
$ws.Cells.Item(78, 10).Value = 10000
$ws.Cells.Item(78, 11).Value = 200
$ws.Cells.Item(78, 12).Value = 200
$ws.Cells.Item(78, 13).Value = 200
$ws.Cells.Item(78, 14).Value = '$/unidad'
$ws.Cells.Item(78, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(78, 16).Value = 200
$ws.Cells.Item(78, 17).Value = 1
$ws.Cells.Item(78, 18).Value = 'Hortaliza'

# Row 79
$ws.Cells.Item(79, 1).Value = 11
$ws.Cells.Item(79, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(79, 3).Value = 'Bíobío'
$ws.Cells.Item(79, 4).Value = 44232
$ws.Cells.Item(79, 5).Value = 8
$ws.Cells.Item(79, 6).Value = 100112024
$ws.Cells.Item(79, 7).Value = 'Choclo'
$ws.Cells.Item(79, 8).Value = 'Choclero'
$ws.Cells.Item(79, 9).Value = 'Primera'
$ws.Cells.Item(79, 10).Value = 10000
$ws.Cells.Item(79, 11).Value = 320
$ws.Cells.Item(79, 12).Value = 350
$ws.Cells.Item(79, 13).Value = 335
$ws.Cells.Item(79, 14).Value = '$/unidad'
$ws.Cells.Item(79, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(79, 16).Value = 335
$ws.Cells.Item(79, 17).Value = 1
$ws.Cells.Item(79, 18).Value = 'Hortaliza'

# Row 80
$ws.Cells.Item(80, 1).Value = 11
$ws.Cells.Item(80, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(80, 3).Value = 'Bíobío'
$ws.Cells.Item(80, 4).Value = 44232
$ws.Cells.Item(80, 5).Value = 8
$ws.Cells.Item(80, 6).Value = 100112024
$ws.Cells.Item(80, 7).Value = 'Choclo'
$ws.Cells.Item(80, 8).Value = 'Choclero'
$ws.Cells.Item(80, 9).Value = 'Segunda'
$ws.Cells.Item(80, 10).Value = 5000
$ws.Cells.Item(80, 11).Value = 280
$ws.Cells.Item(80, 12).Value = 280
$ws.Cells.Item(80, 13).Value = 280
$ws.Cells.Item(80, 14).Value = '$/unidad'
$ws.Cells.Item(80, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(80, 16).Value = 280
$ws.Cells.Item(80, 17).Value = 1
$ws.Cells.Item(80, 18).Value = 'Hortaliza'

# Row 81
$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(81, 3).Value = 'Bíobío'
$ws.Cells.Item(81, 4).Value = 44279
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 100112024
$ws.Cells.Item(81, 7).Value = 'Choclo'
$ws.Cells.Item(81, 8).Value = 'Choclero'
$ws.Cells.Item(81, 9).Value = 'Primera'
$ws.Cells.Item(81, 10).Value = 20000
$ws.Cells.Item(81, 11).Value = 230
$ws.Cells.Item(81, 12).Value = 250
$ws.Cells.Item(81, 13).Value = 240
$ws.Cells.Item(81, 14).Value = '$/unidad'
$ws.Cells.Item(81, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(81, 16).Value = 240
$ws.Cells.Item(81, 17).Value = 1
$ws.Cells.Item(81, 18).Value = 'Hortaliza'

# Row 82
$ws.Cells.Item(82, 1).Value = 11
$ws.Cells.Item(82, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(82, 3).Value = 'Bíobío'
$ws.Cells.Item(82, 4).Value = 44279
$ws.Cells.Item(82, 5).Value = 8
$ws.Cells.Item(82, 6).Value = 100112024
$ws.Cells.Item(82, 7).Value = 'Choclo'
$ws.Cells.Item(82, 8).Value = 'Choclero'
$ws.Cells.Item(82, 9).Value = 'Segunda'
$ws.Cells.Item(82, 10).Value = 5000
$ws.Cells.Item(82, 11).Value = 200
$ws.Cells.Item(82, 12).Value = 200
$ws.Cells.Item(82, 13).Value = 200
$ws.Cells.Item(82, 14).Value = '$/unidad'
$ws.Cells.Item(82, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(82, 16).Value = 200
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = 'Hortaliza'

# Row 83
$ws.Cells.Item(83, 1).Value = 11
$ws.Cells.Item(83, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(83, 3).Value = 'Bíobío'
$ws.Cells.Item(83, 4).Value = 44236
$ws.Cells.Item(83, 5).Value = 8
$ws.Cells.Item(83, 6).Value = 100112024
$ws.Cells.Item(83, 7).Value = 'Choclo'
$ws.Cells.Item(83, 8).Value = 'Choclero'
$ws.Cells.Item(83, 9).Value = 'Primera'
$ws.Cells.Item(83, 10).Value = 10000
$ws.Cells.Item(83, 11).Value = 300
$ws.Cells.Item(83, 12).Value = 320
$ws.Cells.Item(83, 13).Value = 310
$ws.Cells.Item(83, 14).Value = '$/unidad'
$ws.Cells.Item(83, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(83, 16).Value = 310
$ws.Cells.Item(83, 17).Value = 1
$ws.Cells.Item(83, 18).Value = 'Hortaliza'

# Row 84
$ws.Cells.Item(84, 1).Value = 11
$ws.Cells.Item(84, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(84, 3).Value = 'Bíobío'
$ws.Cells.Item(84, 4).Value = 44236
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(84, 6).Value = 100112024
$ws.Cells.Item(84, 7).Value = 'Choclo'
$ws.Cells.Item(84, 8).Value = 'Choclero'
$ws.Cells.Item(84, 9).Value = 'Segunda'
$ws.Cells.Item(84, 10).Value = 5000
$ws.Cells.Item(84, 11).Value = 250
$ws.Cells.Item(84, 12).Value = 250
$ws.Cells.Item(84, 13).Value = 250
$ws.Cells.Item(84, 14).Value = '$/unidad'
$ws.Cells.Item(84, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(84, 16).Value = 250
$ws.Cells.Item(84, 17).Value = 1
$ws.Cells.Item(84, 18).Value = 'Hortaliza'

# Row 85
$ws.Cells.Item(85, 1).Value = 11
$ws.Cells.Item(85, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(85, 3).Value = 'Bíobío'
$ws.Cells.Item(85, 4).Value = 44615
$ws.Cells.Item(85, 5).Value = 8
$ws.Cells.Item(85, 6).Value = 100112024
$ws.Cells.Item(85, 7).Value = 'Choclo'
$ws.Cells.Item(85, 8).Value = 'Choclero'
$ws.Cells.Item(85, 9).Value = 'Primera'
$ws.Cells.Item(85, 10).Value = 10000
$ws.Cells.Item(85, 11).Value = 150
$ws.Cells.Item(85, 12).Value = 160
$ws.Cells.Item(85, 13).Value = 155
$ws.Cells.Item(85, 14).Value = '$/unidad'
$ws.Cells.Item(85, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(85, 16).Value = 155
$ws.Cells.Item(85, 17).Value = 1
$ws.Cells.Item(85, 18).Value = 'Hortaliza'

# Row 86
$ws.Cells.Item(86, 1).Value = 11
$ws.Cells.Item(86, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(86, 3).Value = 'Bíobío'
$ws.Cells.Item(86, 4).Value = 44615
$ws.Cells.Item(86, 5).Value = 8
$ws.Cells.Item(86, 6).Value = 100112024
$ws.Cells.Item(86, 7).Value = 'Choclo'
$ws.Cells.Item(86, 8).Value = 'Choclero'
$ws.Cells.Item(86, 9).Value = 'Segunda'
$ws.Cells.Item(86, 10).Value = 5000
$ws.Cells.Item(86, 11).Value = 120
$ws.Cells.Item(86, 12).Value = 120
$ws.Cells.Item(86, 13).Value = 120
$ws.Cells.Item(86, 14).Value = '$/unidad'
$ws.Cells.Item(86, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(86, 16).Value = 120
$ws.Cells.Item(86, 17).Value = 1
$ws.Cells.Item(86, 18).Value = 'Hortaliza'

# Row 87
$ws.Cells.Item(87, 1).Value = 11
$ws.Cells.Item(87, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(87, 3).Value = 'Bíobío'
$ws.Cells.Item(87, 4).Value = 44349
$ws.Cells.Item(87, 5).Value = 8
$ws.Cells.Item(87, 6).Value = 100112024
$ws.Cells.Item(87, 7).Value = 'Choclo'
$ws.Cells.Item(87, 8).Value = 'Dulce o Americano'
$ws.Cells.Item(87, 9).Value = 'Primera'
$ws.Cells.Item(87, 10).Value = 50
$ws.Cells.Item(87, 11).Value = 9000
$ws.Cells.Item(87, 12).Value = 10000
$ws.Cells.Item(87, 13).Value = 9600
$ws.Cells.Item(87, 14).Value = '$/malla 60 unidades'
$ws.Cells.Item(87, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(87, 16).Value = 160
$ws.Cells.Item(87, 17).Value = 60
$ws.Cells.Item(87, 18).Value = 'Hortaliza'

# Row 88
$ws.Cells.Item(88, 1).Value = 11
$ws.Cells.Item(88, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(88, 3).Value = 'Bíobío'
$ws.Cells.Item(88, 4).Value = 44595
$ws.Cells.Item(88, 5).Value = 8
$ws.Cells.Item(88, 6).Value = 100112024
$ws.Cells.Item(88, 7).Value = 'Choclo'
$ws.Cells.Item(88, 8).Value = 'Choclero'
$ws.Cells.Item(88, 9).Value = 'Primera'
$ws.Cells.Item(88, 10).Value = 20000
$ws.Cells.Item(88, 11).Value = 150
$ws.Cells.Item(88, 12).Value = 200
$ws.Cells.Item(88, 13).Value = 175
$ws.Cells.Item(88, 14).Value = '$/unidad'
$ws.Cells.Item(88, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(88, 16).Value = 175
$ws.Cells.Item(88, 17).Value = 1
$ws.Cells.Item(88, 18).Value = 'Hortaliza'

# Row 89
$ws.Cells.Item(89, 1).Value = 11
$ws.Cells.Item(89, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(89, 3).Value = 'Bíobío'
$ws.Cells.Item(89, 4).Value = 44628
$ws.Cells.Item(89, 5).Value = 8
$ws.Cells.Item(89, 6).Value = 100112024
$ws.Cells.Item(89, 7).Value = 'Choclo'
$ws.Cells.Item(89, 8).Value = 'Choclero'
$ws.Cells.Item(89, 9).Value = 'Primera'
$ws.Cells.Item(89, 10).Value = 3000
$ws.Cells.Item(89, 11).Value = 200
$ws.Cells.Item(89, 12).Value = 200
$ws.Cells.Item(89, 13).Value = 200
$ws.Cells.Item(89, 14).Value = '$/unidad'
$ws.Cells.Item(89, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(89, 16).Value = 200
$ws.Cells.Item(89, 17).Value = 1
$ws.Cells.Item(89, 18).Value = 'Hortaliza'

# Row 90
$ws.Cells.Item(90, 1).Value = 11
$ws.Cells.Item(90, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(90, 3).Value = 'Bíobío'
$ws.Cells.Item(90, 4).Value = 44628
$ws.Cells.Item(90, 5).Value = 8
$ws.Cells.Item(90, 6).Value = 100112024
$ws.Cells.Item(90, 7).Value = 'Choclo'
$ws.Cells.Item(90, 8).Value = 'Choclero'
$ws.Cells.Item(90, 9).Value = 'Segunda'
$ws.Cells.Item(90, 10).Value = 3500
$ws.Cells.Item(90, 11).Value = 150
$ws.Cells.Item(90, 12).Value = 150
$ws.Cells.Item(90, 13).Value = 150
$ws.Cells.Item(90, 14).Value = '$/unidad'
$ws.Cells.Item(90, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(90, 16).Value = 150
$ws.Cells.Item(90, 17).Value = 1
$ws.Cells.Item(90, 18).Value = 'Hortaliza'

# Row 91
$ws.Cells.Item(91, 1).Value = 11
$ws.Cells.Item(91, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(91, 3).Value = 'Bíobío'
$ws.Cells.Item(91, 4).Value = 44552
$ws.Cells.Item(91, 5).Value = 8
$ws.Cells.Item(91, 6).Value = 100112024
$ws.Cells.Item(91, 7).Value = 'Choclo'
$ws.Cells.Item(91, 8).Value = 'Choclero'
$ws.Cells.Item(91, 9).Value = 'Primera'
$ws.Cells.Item(91, 10).Value = 100
$ws.Cells.Item(91, 11).Value = 22000
$ws.Cells.Item(91, 12).Value = 24000
$ws.Cells.Item(91, 13).Value = 23000
$ws.Cells.Item(91, 14).Value = '$/malla 50 unidades'
$ws.Cells.Item(91, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(91, 16).Value = 460
$ws.Cells.Item(91, 17).Value = 50
$ws.Cells.Item(91, 18).Value = 'Hortaliza'

# Row 92
$ws.Cells.Item(92, 1).Value = 11
$ws.Cells.Item(92, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(92, 3).Value = 'Bíobío'
$ws.Cells.Item(92, 4).Value = 44552
$ws.Cells.Item(92, 5).Value = 8
$ws.Cells.Item(92, 6).Value = 100112024
$ws.Cells.Item(92, 7).Value = 'Choclo'
$ws.Cells.Item(92, 8).Value = 'Choclero'
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 100
$ws.Cells.Item(92, 11).Value = 11000
$ws.Cells.Item(92, 12).Value = 12000
$ws.Cells.Item(92, 13).Value = 11500
$ws.Cells.Item(92, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(92, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(92, 16).Value = 164
$ws.Cells.Item(92, 17).Value = 70
$ws.Cells.Item(92, 18).Value = 'Hortaliza'

$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
